# Clean up resultCalibration: drop the no-longer-needed dispersion-distance
# columns (meanDisp Alps / medianDispAlps / meanDispJura / medianDispJura),
# which lived in columns K:N. Deleting them shifts the old "pRepro" values
# column (O) left into K.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("K:N").Delete() | Out-Null

# Leave the selection where the user ended up after the cleanup.
$ws.Range("N11").Select() | Out-Null
